$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure ambiguous numeric-looking text in column D keeps its exact text
# representation (e.g. "1.00", "0.999") instead of being auto-coerced to a
# number by Excel's input parser. We temporarily force Text format on the
# whole D2:D51 data range, write every value, then restore the default
# "Normal" style so the saved XML has no stray style index (matching the
# original unstyled cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.220.86'
$ws.Range("D3").Value = '1.644.91'
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").Value = '216.94'
$ws.Range("D8").Value = '0.258'
$ws.Range("D9").Value = '0.0637'
$ws.Range("D10").Value = '19.81'
$ws.Range("D11").Value = '0.0792'
$ws.Range("D12").Value = '1.874.25'
$ws.Range("D13").Value = '4.28'
$ws.Range("D14").Value = '1.630.36'
$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("D17").Value = '63.25'
$ws.Range("D18").Value = '26.225.52'
$ws.Range("D20").Value = '195.32'
$ws.Range("D21").Value = '4.43'
$ws.Range("D22").Value = '10.06'
$ws.Range("D23").Value = '6.29'
$ws.Range("D25").Value = '1.00'
$ws.Range("D26").Value = '142.87'
$ws.Range("D28").Value = '6.94'
$ws.Range("D29").Value = '15.62'
$ws.Range("D31").Value = '0.0502'
$ws.Range("D32").Value = '3.34'
$ws.Range("D34").Value = '1.59'
$ws.Range("D36").Value = '0.910'
$ws.Range("D37").Value = '1.134.37'
$ws.Range("D38").Value = '0.552'
$ws.Range("D40").Value = '0.0157'
$ws.Range("D42").Value = '100.40'
$ws.Range("D43").Value = '5.50'
$ws.Range("D44").Value = '0.797'
$ws.Range("D45").Value = '1.784.28'
$ws.Range("D46").Value = '56.99'
$ws.Range("D47").Value = '1.48'
$ws.Range("D49").Value = '0.417'
$ws.Range("D50").Value = '7.67'
$ws.Range("D51").Value = '0.0969'

$ws.Range("D2:D51").Style = "Normal"

# Remaining columns (B, C, E) never parse as ambiguous numbers, so a plain
# .Value assignment keeps them as text safely.
$ws.Range("E2").Value = '  +1.73%  '
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("E6").Value = '  +0.98%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("E15").Value = '  -3.04%  '
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("B25").Value = 'BinanceUSD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("E31").Value = '  +1.93%  '
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E50").Value = '  +2.70%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E51").Value = '  +1.93%  '
